$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The s1Protocol (kit) value used for every sample row changed from the old
# NEBNext kit code "E7760" to "E7420". Update column H (s1Protocol) for all
# 16 data rows (rows 2-17); every other column is unchanged.
$ws.Range("H2:H17").Value = "E7420"
